# projekt_useri.xlsx -- "prepravak baze podataka i css"
# Adds a new "id_hotela" column (I) to the user table, backfills it with
# -1 for every existing user row, and appends 12 new user rows (for hotel
# ids 59-70) that carry real id_hotela values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column header -----------------------------------------------
$ws.Cells.Item(1, 9).Value = "id_hotela"

# --- backfill existing rows 2..59 with id_hotela = -1 -----------------
for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 9).Value = -1
}

# --- append new rows 60..71 --------------------------------------------
# Fill column C/D (datum_dolaska / datum_odlaska) with NULL placeholders
# first so the shared-string table grows in the same order as the
# original edit (NULL before the id/username strings).
for ($r = 60; $r -le 71; $r++) {
    $ws.Cells.Item($r, 3).Value = "NULL"
    $ws.Cells.Item($r, 4).Value = "NULL"
}

# id_usera values for the new rows -- stored as text, like the rest of
# column A, even though they look numeric.
$ids = @("59", "60", "61", "62", "63", "64", "65", "66", "67", "68", "69", "70")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $cell = $ws.Cells.Item(60 + $i, 1)
    $cell.Value = "'" + $ids[$i]
    $cell.Style = "Normal"
}

# usernames for the new rows
$names = @("Ivana1", "Kiki1", "Dorotea1", "Ivana2", "Kiki2", "Dorotea2", "Ivana3", "Kiki3", "Dorotea3", "Ivana4", "Kiki4", "Dorotea4")
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item(60 + $i, 2).Value = $names[$i]
}

# password_hash (E), email (F), registration_sequence (G), has_registered
# (H) and the new id_hotela (I) values for every new row.
$eVals = @(181, 182, 183, 184, 185, 186, 187, 188, 189, 190, 191, 192)
$iVals = @(1, 2, 3, 6, 7, 11, 12, 13, 16, 17, 21, 22)
for ($i = 0; $i -lt 12; $i++) {
    $r = 60 + $i
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
    $ws.Cells.Item($r, 6).Value = "a@b.com"
    $ws.Cells.Item($r, 7).Value = "abc"
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 9).Value = $iVals[$i]
}

# --- restore the selection left on the sheet after the edit -----------
[void]$ws.Range("H61").Select()
